$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'261.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.86%"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'27.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.05%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'4.706"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.44%"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.06204"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.68%"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'6.719"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.67%"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.8499"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.99%"
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.9167"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-1.91%"
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.1408"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.81%"
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.04614"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-8.73%"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07084"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.39%"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.03152"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.73%"
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.09056"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.90%"
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.001532"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.85%"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.0006157"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.34%"
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.006025"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.56%"
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'3.468"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.09%"
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'3.167"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.24%"
$ws.Range("E18").Style = "Normal"

$ws.Range("E20").Value = "'0.40%"
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Value = "'0.99%"
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'4.102"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-0.32%"
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.04222"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.12%"
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.001212"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.07%"
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = "'-6.02%"
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.0001201"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'0.31%"
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'0.0001601"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'5.32%"
$ws.Range("E27").Style = "Normal"

$ws.Range("D40").Value = "'0.03925"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.60%"
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'-0.31%"
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.004132"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'4.93%"
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = "'-4.62%"
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.01379"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-9.74%"
$ws.Range("E44").Style = "Normal"

$ws.Range("E45").Value = "'0.98%"
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Value = "'0.27%"
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'-34.05%"
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "'27.44%"
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.27%"
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.0002001"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.27%"
$ws.Range("E50").Style = "Normal"
